$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.350.64'
$ws.Range("E2").Value = '  +0.40%  '
$ws.Range("D3").Value = '1.869.08'
$ws.Range("D4").Value = '''1.001'
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").Value = '''237.89'
$ws.Range("E5").Value = '  +0.67%  '
$ws.Range("D6").Value = '''1.001'
$ws.Range("E6").Value = '  +0.06%  '
$ws.Range("D7").Value = '''0.4815'
$ws.Range("E7").Value = '  -0.48%  '
$ws.Range("D8").Value = '''0.2806'
$ws.Range("E8").Value = '  -2.24%  '
$ws.Range("D9").Value = '''0.06510'
$ws.Range("E9").Value = '  -1.09%  '
$ws.Range("D10").Value = '1.865.27'
$ws.Range("E10").Value = '  -0.74%  '
$ws.Range("D11").Value = '''0.07463'
$ws.Range("E11").Value = '  +1.92%  '
$ws.Range("D12").Value = '''16.52'
$ws.Range("E12").Value = '  -1.50%  '
$ws.Range("D13").Value = '''5.078'
$ws.Range("E13").Value = '  -0.95%  '
$ws.Range("D14").Value = '''87.81'
$ws.Range("E14").Value = '  +0.86%  '
$ws.Range("D15").Value = '''0.6550'
$ws.Range("E15").Value = '  +0.10%  '
$ws.Range("D16").Value = '30.302.11'
$ws.Range("E16").Value = '  +0.32%  '
$ws.Range("D17").Value = '''13.26'
$ws.Range("E17").Value = '  -0.45%  '
$ws.Range("E18").Value = '  +0.20%  '
$ws.Range("D19").Value = '''0.000007594'
$ws.Range("E19").Value = '  -1.94%  '
$ws.Range("D20").Value = '2.107.06'
$ws.Range("E20").Value = '  -0.99%  '
$ws.Range("D21").Value = '''5.292'
$ws.Range("E21").Value = '  -1.72%  '
$ws.Range("D22").Value = '''1.001'
$ws.Range("E22").Value = '  +0.14%  '
$ws.Range("D23").Value = '''220.50'
$ws.Range("E23").Value = '  +14.36%  '
$ws.Range("D24").Value = '''6.163'
$ws.Range("E24").Value = '  +0.67%  '
$ws.Range("D25").Value = '''9.304'
$ws.Range("E25").Value = '  +0.55%  '
$ws.Range("D26").Value = '''167.46'
$ws.Range("E26").Value = '  +2.88%  '
$ws.Range("D27").Value = '''18.51'
$ws.Range("E27").Value = '  +2.84%  '
$ws.Range("D28").Value = '''1.966'
$ws.Range("E28").Value = '  +2.95%  '
$ws.Range("D29").Value = '''1.450'
$ws.Range("E29").Value = '  +1.17%  '
$ws.Range("D30").Value = '''0.09341'
$ws.Range("E30").Value = '  +2.92%  '
$ws.Range("D31").Value = '''4.313'
$ws.Range("E31").Value = '  +1.30%  '
$ws.Range("D32").Value = '''4.027'
$ws.Range("E32").Value = '  +0.51%  '
$ws.Range("D33").Value = '''0.05054'
$ws.Range("E33").Value = '  -0.05%  '
$ws.Range("D34").Value = '''1.208'
$ws.Range("E34").Value = '  +10.06%  '
$ws.Range("D35").Value = '''0.7479'
$ws.Range("E35").Value = '  +4.65%  '
$ws.Range("D36").Value = '''2.716'
$ws.Range("E36").Value = '  +0.60%  '
$ws.Range("D37").Value = '''0.01831'
$ws.Range("E37").Value = '  +3.04%  '
$ws.Range("D38").Value = '''2.616'
$ws.Range("E38").Value = '  -0.77%  '
$ws.Range("B39").Value = 'RenderToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D39").Value = '''2.086'
$ws.Range("E39").Value = '  +2.10%  '
$ws.Range("B40").Value = 'TrustWalletToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D40").Value = '''0.9121'
$ws.Range("E40").Value = '  -1.13%  '
$ws.Range("D41").Value = '''106.92'
$ws.Range("E41").Value = '  +1.11%  '
$ws.Range("D42").Value = '''5.921'
$ws.Range("E42").Value = '  +2.35%  '
$ws.Range("D43").Value = '''0.4270'
$ws.Range("E43").Value = '  +0.07%  '
$ws.Range("E44").Value = '  +0.40%  '
$ws.Range("D45").Value = '''7.414'
$ws.Range("E45").Value = '  +0.47%  '
$ws.Range("D46").Value = '''0.1293'
$ws.Range("E46").Value = '  -1.09%  '
$ws.Range("D47").Value = '''64.09'
$ws.Range("E47").Value = '  -1.01%  '
$ws.Range("B48").Value = 'NEARProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D48").Value = '''1.478'
$ws.Range("E48").Value = '  +8.72%  '
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").Value = '''8.955'
$ws.Range("E49").Value = '  +1.48%  '
$ws.Range("D50").Value = '''33.75'
$ws.Range("E50").Value = '  +0.16%  '
$ws.Range("B51").Value = 'Decentraland'
$ws.Range("C51").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D51").Value = '''0.3880'
$ws.Range("E51").Value = '  +1.86%  '
